$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows before row 194, shifting the existing
# rows 194-258 down to 196-260.
$ws.Rows.Item(194).Resize(2).Insert()

# Populate the two newly inserted rows (194 and 195) with their data.
$ws.Cells.Item(194, 1).Value = 9
$ws.Cells.Item(194, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(194, 3).Value = "Metropolitana"
$ws.Cells.Item(194, 4).Value = 44468
$ws.Cells.Item(194, 5).Value = 13
$ws.Cells.Item(194, 6).Value = 100112039
$ws.Cells.Item(194, 7).Value = "Ciboulette"
$ws.Cells.Item(194, 8).Value = "Sin especificar"
$ws.Cells.Item(194, 9).Value = "Primera"
$ws.Cells.Item(194, 10).Value = 250
$ws.Cells.Item(194, 11).Value = 1000
$ws.Cells.Item(194, 12).Value = 1200
$ws.Cells.Item(194, 13).Value = 1100
$ws.Cells.Item(194, 14).Value = "$/docena de atados"
$ws.Cells.Item(194, 15).Value = "Región Metropolitana"
$ws.Cells.Item(194, 16).Value = 367
$ws.Cells.Item(194, 17).Value = 3
$ws.Cells.Item(194, 18).Value = "Hortaliza"

$ws.Cells.Item(195, 1).Value = 9
$ws.Cells.Item(195, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(195, 3).Value = "Metropolitana"
$ws.Cells.Item(195, 4).Value = 44468
$ws.Cells.Item(195, 5).Value = 13
$ws.Cells.Item(195, 6).Value = 100112039
$ws.Cells.Item(195, 7).Value = "Ciboulette"
$ws.Cells.Item(195, 8).Value = "Sin especificar"
$ws.Cells.Item(195, 9).Value = "Segunda"
$ws.Cells.Item(195, 10).Value = 131
$ws.Cells.Item(195, 11).Value = 700
$ws.Cells.Item(195, 12).Value = 900
$ws.Cells.Item(195, 13).Value = 801
$ws.Cells.Item(195, 14).Value = "$/docena de atados"
$ws.Cells.Item(195, 15).Value = "Región Metropolitana"
$ws.Cells.Item(195, 16).Value = 267
$ws.Cells.Item(195, 17).Value = 3
$ws.Cells.Item(195, 18).Value = "Hortaliza"

# Match the date style (numFmt) used by the rest of column D.
$ws.Cells.Item(194, 4).NumberFormat = $ws.Cells.Item(196, 4).NumberFormat
$ws.Cells.Item(195, 4).NumberFormat = $ws.Cells.Item(196, 4).NumberFormat

Write-Output "done"
